# Aggiornamento dati fino al 9 settembre 2021 (compreso)
# Adds rows 367-374 (dates 2021-09-02 .. 2021-09-09) to the data table,
# mirroring the formatting (date number-format) already used in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 367; Date = 44441; B = 2; C = 3; D = 122.3990208078335 },
    @{ Row = 368; Date = 44442; B = 1; C = 4; D = 163.1986944104447 },
    @{ Row = 369; Date = 44443; B = 0; C = 3; D = 122.3990208078335 },
    @{ Row = 370; Date = 44444; B = 0; C = 3; D = 122.3990208078335 },
    @{ Row = 371; Date = 44445; B = 1; C = 4; D = 163.1986944104447 },
    @{ Row = 372; Date = 44446; B = 0; C = 4; D = 163.1986944104447 },
    @{ Row = 373; Date = 44447; B = 0; C = 4; D = 163.1986944104447 },
    @{ Row = 374; Date = 44448; B = 0; C = 2; D = 81.59934720522236 }
)

foreach ($r in $newRows) {
    # Copy column-A formatting (the date number format applied to A2:A366)
    # down onto the new date cell before writing its value.
    $ws.Range("A366").Copy()
    $ws.Range("A" + $r.Row).PasteSpecial(-4122)

    $ws.Range("A" + $r.Row).Value2 = $r.Date
    $ws.Range("B" + $r.Row).Value2 = $r.B
    $ws.Range("C" + $r.Row).Value2 = $r.C
    $ws.Range("D" + $r.Row).Value2 = $r.D
}
